$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 37.61313626083334
$ws.Range("C3").Value = 38.43518730115238
$ws.Range("C4").Value = 39.66689682492857
$ws.Range("C5").Value = 38.41266454308918
$ws.Range("C6").Value = 34.99892058125839
